$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-116
# from serial date 45186 (2023-09-17) to 45188 (2023-09-19)
$ws.Range("C2:C116").Value = 45188
